$d = $word.ActiveDocument

$replacements = @(
    @{old = "98×84=8232"; new = "49×88=4312"},
    @{old = "21×97=2037"; new = "58×18=1044"},
    @{old = "91×97=8827"; new = "48×93=4464"},
    @{old = "61×27=1647"; new = "66×25=1650"},
    @{old = "27×19=513";  new = "45×17=765"},
    @{old = "60×87=5220"; new = "47×56=2632"},
    @{old = "20×95=1900"; new = "96×92=8832"},
    @{old = "93×13=1209"; new = "26×62=1612"},
    @{old = "67×22=1474"; new = "95×82=7790"},
    @{old = "19×70=1330"; new = "59×34=2006"},
    @{old = "60×37=2220"; new = "19×67=1273"},
    @{old = "26×27=702";  new = "27×51=1377"},
    @{old = "28×39=1092"; new = "57×49=2793"},
    @{old = "64×29=1856"; new = "26×62=1612"},
    @{old = "55×79=4345"; new = "12×44=528"},
    @{old = "18×19=342";  new = "24×57=1368"},
    @{old = "17×68=1156"; new = "43×91=3913"},
    @{old = "14×22=308";  new = "98×96=9408"},
    @{old = "94×96=9024"; new = "77×97=7469"},
    @{old = "43×40=1720"; new = "40×98=3920"},
    @{old = "12×31=372";  new = "45×63=2835"},
    @{old = "67×47=3149"; new = "48×96=4608"},
    @{old = "67×52=3484"; new = "19×20=380"},
    @{old = "41×77=3157"; new = "38×67=2546"},
    @{old = "12×51=612";  new = "90×84=7560"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
